$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text change: "Main Shelves" -> "Main Shelf" ---
# All cells A9:A15 share the same string ("Main Shelves"); rewrite them all so
# the shared-string table collapses the edit in place instead of branching a
# brand new entry.
$ws.Range("A9:A15").Value = "Main Shelf"

# --- Selection change: B19 -> A9 ---
$ws.Range("A9").Select()

# --- Column width changes ---
$ws.Columns.Item(1).ColumnWidth = 27.833333333333332
$ws.Columns.Item(2).ColumnWidth = 56.833333333333336
$ws.Range("C:D").ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 10.5
$ws.Columns.Item(6).ColumnWidth = 64.16666666666667
$ws.Range("G1:AMK1").EntireColumn.ColumnWidth = 10.5
